$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "RHPF"
# ---------------------------------------------------------------------
$rhpf = $wb.Worksheets.Item("RHPF")

# Two new pathway columns / rows are being added to the matrix
$rhpf.Range("G1").Value = "electrolysis with guaranteed clean electricity"
$rhpf.Range("H1").Value = "natural gas reforming with CCS"
$rhpf.Range("G1:H1").WrapText = $true

$rhpf.Range("A7").Value = "electrolysis with guaranteed clean electricity"
$rhpf.Range("A8").Value = "natural gas reforming with CCS"

# Data block used to route everything to "electrolysis"; it now routes
# everything to "electrolysis with guaranteed clean electricity"
$rhpf.Range("B2:H6").Value = 0
$rhpf.Range("B7:H7").Value = 1
$rhpf.Range("B8:H8").Value = 0

$rhpf.Columns("G:H").ColumnWidth = 16.5

$rhpf.Rows("1:1").Select()

# ---------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# A7:A9 lose their (unused) bold-font style -> back to default formatting
$about.Range("A7:A9").ClearFormats()

# A12 used to read "electrolysis." - now a fuller description
$about.Range("A12").Value = "electrolysis that is guaranteed to be supplied by new clean electricity sources."

$about.Range("A13").Select()
